$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The pH reading for "31 ppm P" / week 02 was originally recorded as the text
# "6.87*" (flagged because the pot was overfilled). After correcting for the
# overfill, the real measured pH is 6.96, so replace the flagged text value
# with the corrected numeric reading.
$ws.Range("C8").Value = 6.96

# Update the Notes cell to the right of "Notes" with the expanded note that
# explains the overfill, the original (flagged) reading, and the corrected
# follow-up reading/date.
$ws.Range("C12").Value = "*Overfilled by 675 mL on 08/15/2025. pH was 6.87, corrected by adding extra nut. On 08/18/2025, new pH. "

# The "Notes" label cell had picked up an empty/no-op alignment style -
# clear its formatting so it goes back to the default (unstyled) cell.
$ws.Range("A12").ClearFormats()

# Row 13 was a stray empty styled row below the notes - remove it entirely.
$ws.Rows(13).Delete()

# Update the active selection to match the saved view state.
$ws.Range("K7").Select()
